$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.931.21"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.554.33"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.83"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.91"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "1.775.70"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "1.555.16"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("D16").Value = "26.924.73"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.68"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.83"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.21"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.92"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0467"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "1.436.84"
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.973"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0165"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.69"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.02"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").Value = "1.689.37"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.46"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0525"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "0.0₇0980"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  +1.60%  "
